$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row groups whose entire contents get cyclically rotated:
# value of row cyc[i] becomes the old value of row cyc[i+1] (wrapping).
# (Derived from the target diff: full-row data is re-sorted while the
# row indices/number of rows on the sheet stay fixed.)
$cycles = @(
    ,@(2, 3)
    ,@(7, 8, 9)
    ,@(10, 11)
    ,@(18, 21, 20, 19)
    ,@(22, 23)
    ,@(27, 28)
    ,@(34, 35, 36)
    ,@(45, 46)
)

# Columns A..AY (1..51), skipping Y (25) and AA (27): those two hold the
# "Startdatum"/"Slutdatum" text values, which are identical across every
# row inside each cycle, so they never need to move; skipping them also
# avoids Excel's automatic text->date conversion on reassignment.
$cols = @()
for ($c = 1; $c -le 51; $c++) {
    if ($c -ne 25 -and $c -ne 27) {
        $cols += $c
    }
}

foreach ($cycle in $cycles) {
    $n = $cycle.Length

    # Snapshot the current value of every cell in every row of this cycle.
    $snapshot = @{}
    foreach ($row in $cycle) {
        $rowVals = @{}
        foreach ($c in $cols) {
            $rowVals[$c] = $ws.Cells.Item($row, $c).Value()
        }
        $snapshot[$row] = $rowVals
    }

    # Write back rotated: row[i] <- row[i+1] (wrapping around).
    for ($i = 0; $i -lt $n; $i++) {
        $dstRow = $cycle[$i]
        $srcRow = $cycle[($i + 1) % $n]
        $srcVals = $snapshot[$srcRow]
        foreach ($c in $cols) {
            $ws.Cells.Item($dstRow, $c).Value = $srcVals[$c]
        }
    }
}
